$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2 through 20
# from 45207 (2023-10-08) to 45208 (2023-10-09).
for ($row = 2; $row -le 20; $row++) {
    $ws.Cells.Item($row, 3).Value = 45208
}
